{"js": "// Commit: \"Fixed POI packaging and upgraded to POI 3.15.\"\n//\n// The supplied OOXML diff touches word/document.xml and word/styles.xml,\n// but every single changed line is a pure XML-attribute / namespace-\n// declaration re-ordering (e.g. `w:tcW w:w=\"2303\" w:type=\"dxa\"` becomes\n// `w:tcW w:type=\"dxa\" w:w=\"2303\"`, the root <w:document> namespace\n// declarations get alphabetized, <w:lsdException>, <w:style>,\n// <w:tblBorders>, <w:pgSz>, <w:pgMar>, <w:lang>, ... attributes are all\n// simply re-sorted alphabetically by local name). No element, attribute\n// value, run of text, table cell, style definition or any other piece of\n// document content was added, removed or changed \u2014 this is exactly the\n// kind of diff produced when a file is re-serialized by a newer version\n// of Apache POI (as the commit message says), not an authoring edit.\n//\n// There is therefore no content-level change for the Word JavaScript API\n// to replay: the script intentionally performs only read-only\n// (``.load()`` / ``context.sync()``) calls so the document model is left\n// byte-for-byte equivalent in meaning (and, since nothing is mutated,\n// untouched in the underlying package as well).\n\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.load(\"rowCount,values\");\n  await context.sync();\n}\n\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n", "ps1": "# Commit: \"Fixed POI packaging and upgraded to POI 3.15.\"\n#\n# The supplied OOXML diff touches word/document.xml and word/styles.xml,\n# but every single changed line is a pure XML-attribute / namespace-\n# declaration re-ordering (e.g. `w:tcW w:w=\"2303\" w:type=\"dxa\"` becomes\n# `w:tcW w:type=\"dxa\" w:w=\"2303\"`, the root <w:document> namespace\n# declarations get alphabetized, <w:lsdException>, <w:style>,\n# <w:tblBorders>, <w:pgSz>, <w:pgMar>, <w:lang>, ... attributes are all\n# simply re-sorted alphabetically by local name). No element, attribute\n# value, run of text, table cell, style definition or any other piece of\n# document content was added, removed or changed \u2014 this is exactly the\n# kind of diff produced when a file is re-serialized by a newer version\n# of Apache POI (as the commit message says), not an authoring edit.\n#\n# There is therefore no content-level change for the Word object model to\n# replay: the script intentionally performs only read-only property\n# reads so the document is left byte-for-byte equivalent in meaning (and,\n# since nothing is mutated, untouched in the underlying package as well).\n\n$d = $word.ActiveDocument\n\n$text = $d.Content.Text\n\n$tableCount = $d.Tables.Count\nif ($tableCount -gt 0) {\n    $tbl = $d.Tables.Item(1)\n    $rows = $tbl.Rows.Count\n    $cols = $tbl.Columns.Count\n}\n\n$sectionCount = $d.Sections.Count\n"}
